$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency market-data snapshot table (prices & ranking reshuffle)

$ws.Range("D2").Value = "'243.10"

$ws.Range("D3").Value = "'23.84"

$ws.Range("D4").Value = "'5.238"

$ws.Range("D5").Value = "'0.05778"

$ws.Range("D6").Value = "'6.437"

$ws.Range("D7").Value = "'3.235"

$ws.Range("D8").Value = "'0.8062"

$ws.Range("D9").Value = "'0.8799"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1388"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07101"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03106"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03056"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09314"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.819"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001529"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04710"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0006043"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").Value = "'0.006161"

$ws.Range("D20").Value = "'0.001262"

$ws.Range("D21").Value = "'0.004061"

$ws.Range("D22").Value = "'0.00008723"

$ws.Range("D23").Value = "'3.543"

$ws.Range("D24").Value = "'2.156"

$ws.Range("D26").Value = "'0.1320"

$ws.Range("D28").Value = "'0.0002333"

$ws.Range("D40").Value = "'0.03748"

$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1050"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002519"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003247"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

$ws.Range("D44").Value = "'0.007831"

$ws.Range("D45").Value = "'0.00005334"

$ws.Range("D46").Value = "'0.00000000752"

$ws.Range("D47").Value = "'0.5361"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D48").Value = "'0.002496"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"

$ws.Range("D49").Value = "'0.00002104"

$ws.Range("D50").Value = "'0.0002004"
